$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.158.89'
$ws.Range("E2").Value = '  -0.25%  '
$ws.Range("D3").Value = '1.584.18'
$ws.Range("E3").Value = '  -0.07%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("E5").Value = '  +0.87%  '
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = '0.0606'
$ws.Range("E9").Value = '  -0.95%  '
$ws.Range("E10").Value = '  -1.89%  '
$ws.Range("E11").Value = '  +0.26%  '
$ws.Range("D12").Value = '1.806.87'
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("D13").Value = '1.594.14'
$ws.Range("E13").Value = '  +0.05%  '
$ws.Range("E14").Value = '  -1.54%  '
$ws.Range("E15").Value = '  -0.04%  '
$ws.Range("D16").Value = '64.05'
$ws.Range("E16").Value = '  -0.77%  '
$ws.Range("D17").Value = '26.174.05'
$ws.Range("E17").Value = '  -0.15%  '
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("E19").Value = '  -0.72%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '212.80'
$ws.Range("E20").Value = '  +0.55%  '
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("E22").Value = '  -0.60%  '
$ws.Range("D23").Value = '2.17'
$ws.Range("E23").Value = '  -0.67%  '
$ws.Range("D24").Value = '8.95'
$ws.Range("E24").Value = '  +1.08%  '
$ws.Range("D25").Value = '143.62'
$ws.Range("E25").Value = '  -0.57%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").Value = '6.96'
$ws.Range("E27").Value = '  -0.88%  '
$ws.Range("E28").Value = '  -0.86%  '
$ws.Range("D29").Value = '15.11'
$ws.Range("E29").Value = '  -1.29%  '
$ws.Range("D30").Value = '0.0496'
$ws.Range("E30").Value = '  -2.34%  '
$ws.Range("E31").Value = '  +0.50%  '
$ws.Range("E32").Value = '  -1.29%  '
$ws.Range("D33").Value = '1.339.21'
$ws.Range("E33").Value = '  +4.09%  '
$ws.Range("E34").Value = '  -2.11%  '
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("E36").Value = '  -1.54%  '
$ws.Range("E37").Value = '  -3.81%  '
$ws.Range("E38").Value = '  -0.52%  '
$ws.Range("E39").Value = '  +0.25%  '
$ws.Range("E40").Value = '  +3.26%  '
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("E42").Value = '  -15.12%  '
$ws.Range("E43").Value = '  +0.55%  '
$ws.Range("E44").Value = '  -0.16%  '
$ws.Range("D45").Value = '1.719.83'
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("D46").Value = '60.91'
$ws.Range("E46").Value = '  -2.58%  '
$ws.Range("D47").Value = '85.88'
$ws.Range("E47").Value = '  -3.28%  '
$ws.Range("D48").Value = '0.0₆0103'
$ws.Range("E48").Value = '  +6.29%  '
$ws.Range("E49").Value = '  -2.02%  '
$ws.Range("D50").Value = '0.0984'
$ws.Range("E50").Value = '  -1.89%  '
$ws.Range("E51").Value = '  -1.02%  '
